$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Formatting: header / title styling -----------------------------------
# The title (A1, merged A1:G1) loses its 14pt size (falls back to the
# workbook default of 11pt) and becomes white so it reads clearly against
# its background; the column headers on row 2 (already bold, on the blue
# fill) also become white so they're legible on the dark blue header band.
$title = $ws.Range("A1")
$title.Font.Size = 11
$title.Font.Color = 16777215   # RGB(255,255,255) -> white

$headerRow = $ws.Range("A2:K2")
$headerRow.Font.Color = 16777215   # RGB(255,255,255) -> white

# --- Data updates on row 3 --------------------------------------------------
# PERIOD TO EXPIRE drops from 126 to 118 days.
$ws.Range("H3").Value = 118

# LAST UPDATE moves from 08-Sep-2025 to 16-Sep-2025. This column stores its
# dates as plain text, so we stage the new text through a scratch cell
# that's pre-formatted as Text and paste only the value in, which keeps
# Excel from reinterpreting the string as a serial date.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "16-Sep-2025"
$scratch.Copy()
$ws.Range("I3").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()
